$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same date serial (45178 = 2023-09-09) for
# every data row (rows 2-238). The update bumps it by one day to
# 45179 (2023-09-10) across the whole column.
$ws.Range("C2:C238").Value = 45179
